$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Double the values in A2:A10 (A1 stays unchanged at 0)
for ($i = 2; $i -le 10; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $current = $cell.Value()
    $cell.Value = $current * 2
}

# Update the chart: title text and the plotted data range
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "LINE CHART"

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,,'Sheet'!`$A`$1:`$A`$5,1)"
